# Refresh the cryptocurrency price / 1h-volume table with the latest scrape.
# Two pairs of adjacent rows (37/38, 42/43, 50/51) also had their coin
# identity (Coin name + Link) swap rank position, per the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.950.10"
$ws.Range("E2").Value = "  -3.59%  "
$ws.Range("D3").Value = "1.869.26"
$ws.Range("E3").Value = "  -2.46%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9996"
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.00"
$ws.Range("E5").Value = "  -2.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4369"
$ws.Range("E7").Value = "  -4.98%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3720"
$ws.Range("E8").Value = "  -2.85%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07501"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9385"
$ws.Range("E10").Value = "  -4.37%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.38"
$ws.Range("E11").Value = "  -3.81%  "
$ws.Range("D12").Value = "1.877.24"
$ws.Range("E12").Value = "  -3.15%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.757"
$ws.Range("E13").Value = "  -2.93%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.460"
$ws.Range("E14").Value = "  -4.08%  "
$ws.Range("E15").Value = "  -1.79%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("E16").Value = "  -0.16%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "81.73"
$ws.Range("E17").Value = "  -2.79%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009095"
$ws.Range("E18").Value = "  -4.14%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.001"
$ws.Range("E19").Value = "  +0.00%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.97"
$ws.Range("E20").Value = "  -4.26%  "
$ws.Range("D21").Value = "27.937.20"
$ws.Range("E21").Value = "  -3.63%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.132"
$ws.Range("E22").Value = "  -3.62%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.05"
$ws.Range("E23").Value = "  +0.95%  "
$ws.Range("D24").Value = "2.119.48"
$ws.Range("E24").Value = "  -3.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.002"
$ws.Range("E25").Value = "  -3.92%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.28"
$ws.Range("E26").Value = "  -2.67%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.51"
$ws.Range("E27").Value = "  -3.09%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.506"
$ws.Range("E28").Value = "  -3.28%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "113.37"
$ws.Range("E29").Value = "  -3.64%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.727"
$ws.Range("E30").Value = "  -7.69%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09039"
$ws.Range("E31").Value = "  -2.97%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8238"
$ws.Range("E32").Value = "  -4.92%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.829"
$ws.Range("E33").Value = "  -5.44%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.179"
$ws.Range("E34").Value = "  -5.69%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.944"
$ws.Range("E35").Value = "  -3.52%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.001"
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05510"
$ws.Range("E37").Value = "  -3.51%  "
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.121"
$ws.Range("E38").Value = "  -3.14%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01984"
$ws.Range("E39").Value = "  -3.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.946"
$ws.Range("E40").Value = "  -3.13%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5276"
$ws.Range("E41").Value = "  -4.38%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.066"
$ws.Range("E42").Value = "  -5.94%  "
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1710"
$ws.Range("E43").Value = "  -2.47%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.818"
$ws.Range("E44").Value = "  -6.12%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.06767"
$ws.Range("E45").Value = "  -1.87%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4921"
$ws.Range("E46").Value = "  -5.17%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.69"
$ws.Range("E47").Value = "  -4.45%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "107.55"
$ws.Range("E48").Value = "  -2.51%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.685"
$ws.Range("E49").Value = "  -5.55%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.907"
$ws.Range("E50").Value = "  -12.58%  "
$ws.Range("B51").Value = "PaxDollar"
$ws.Range("C51").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9993"
$ws.Range("E51").Value = "  -0.11%  "
